# Auto-generated cell updates derived from the OOXML diff.
$wb = $excel.ActiveWorkbook

function Set-GoblinCell {
    param($ws, $row, $col, $value)
    $ws.Cells.Item($row, $col).Value = $value
}

function Clear-GoblinCell {
    param($ws, $row, $col)
    $ws.Cells.Item($row, $col).ClearContents()
}

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-GoblinCell $ws 64 8 8581.444   # H64: 8403.536 -> 8581.444
Set-GoblinCell $ws 64 9 3959.8   # I64: 4079.8 -> 3959.8
Set-GoblinCell $ws 64 10 9631.817999999999   # J64: 9343.478999999999 -> 9631.817999999999
Set-GoblinCell $ws 64 11 3959.8   # K64: 4079.8 -> 3959.8
Set-GoblinCell $ws 64 12 9631.817999999999   # L64: 9343.478999999999 -> 9631.817999999999
Set-GoblinCell $ws 64 13 -3711.8   # M64: -3831.8 -> -3711.8
Set-GoblinCell $ws 64 14 -10127.818   # N64: -9839.478999999999 -> -10127.818
Set-GoblinCell $ws 67 8 8581.444   # H67: 8403.536 -> 8581.444
Set-GoblinCell $ws 67 9 3959.8   # I67: 4079.8 -> 3959.8
Set-GoblinCell $ws 67 10 9631.817999999999   # J67: 9343.478999999999 -> 9631.817999999999
Set-GoblinCell $ws 67 11 3959.8   # K67: 4079.8 -> 3959.8
Set-GoblinCell $ws 67 12 9631.817999999999   # L67: 9343.478999999999 -> 9631.817999999999
Set-GoblinCell $ws 67 13 -3101.8   # M67: -3221.8 -> -3101.8
Set-GoblinCell $ws 67 14 -11347.818   # N67: -11059.479 -> -11347.818
Set-GoblinCell $ws 98 8 1727.8518   # H98: 1870.4445 -> 1727.8518
Set-GoblinCell $ws 98 9 1466.36   # I98: 1521.2084 -> 1466.36
Set-GoblinCell $ws 98 10 4996.5   # J98: 4664.3335 -> 4996.5
Set-GoblinCell $ws 98 11 1466.36   # K98: 1521.2084 -> 1466.36
Set-GoblinCell $ws 98 12 4996.5   # L98: 4664.3335 -> 4996.5
Set-GoblinCell $ws 98 13 31.6400000000001   # M98: -23.20839999999998 -> 31.6400000000001
Set-GoblinCell $ws 98 14 -7992.5   # N98: -7660.3335 -> -7992.5
Set-GoblinCell $ws 122 8 1727.8518   # H122: 1870.4445 -> 1727.8518
Set-GoblinCell $ws 122 9 1466.36   # I122: 1521.2084 -> 1466.36
Set-GoblinCell $ws 122 10 4996.5   # J122: 4664.3335 -> 4996.5
Set-GoblinCell $ws 122 11 4399.08   # K122: 4563.6252 -> 4399.08
Set-GoblinCell $ws 122 12 14989.5   # L122: 13993.0005 -> 14989.5
Set-GoblinCell $ws 122 13 -1949.08   # M122: -2113.6252 -> -1949.08
Set-GoblinCell $ws 122 14 -19889.5   # N122: -18893.0005 -> -19889.5
Set-GoblinCell $ws 129 8 1194.1666   # H129: 1194.4166 -> 1194.1666
Set-GoblinCell $ws 129 10 2184.2222   # J129: 2184.889 -> 2184.2222
Set-GoblinCell $ws 129 12 6552.6666   # L129: 6554.667 -> 6552.6666
Set-GoblinCell $ws 129 14 -16552.6666   # N129: -16554.667 -> -16552.6666
Set-GoblinCell $ws 132 8 1995.1316   # H132: 2049.6487 -> 1995.1316
Set-GoblinCell $ws 132 9 1435.2059   # I132: 1479.3636 -> 1435.2059
Set-GoblinCell $ws 132 11 4305.6177   # K132: 4438.0908 -> 4305.6177
Set-GoblinCell $ws 132 13 -1775.6177   # M132: -1908.0908 -> -1775.6177
Set-GoblinCell $ws 135 8 1377.3684   # H135: 1457.1666 -> 1377.3684
Set-GoblinCell $ws 135 9 1288.75   # I135: 1329.9375 -> 1288.75
Set-GoblinCell $ws 135 10 1850   # J135: 2475 -> 1850
Set-GoblinCell $ws 135 11 11598.75   # K135: 11969.4375 -> 11598.75
Set-GoblinCell $ws 135 12 16650   # L135: 22275 -> 16650
Set-GoblinCell $ws 135 13 -9063.75   # M135: -9434.4375 -> -9063.75
Set-GoblinCell $ws 135 14 -21720   # N135: -27345 -> -21720
Set-GoblinCell $ws 138 8 3505.9158   # H138: 3507.0435 -> 3505.9158
Set-GoblinCell $ws 138 10 4064.4265   # J138: 4091.8 -> 4064.4265
Set-GoblinCell $ws 138 12 12193.2795   # L138: 12275.4 -> 12193.2795
Set-GoblinCell $ws 138 14 -22473.2795   # N138: -22555.4 -> -22473.2795

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-GoblinCell $ws 3 8 4500   # H3: 2750 -> 4500
Set-GoblinCell $ws 3 9 4500   # I3: 2750 -> 4500
Set-GoblinCell $ws 3 11 4500   # K3: 2750 -> 4500
Set-GoblinCell $ws 3 13 -4385   # M3: -2635 -> -4385
Set-GoblinCell $ws 11 8 4201002   # H11: 6835001 -> 4201002
Set-GoblinCell $ws 11 9 5250003   # I11: 10250002 -> 5250003
Set-GoblinCell $ws 11 10 4999   # J11: 5000 -> 4999
Set-GoblinCell $ws 11 11 5250003   # K11: 10250002 -> 5250003
Set-GoblinCell $ws 11 12 4999   # L11: 5000 -> 4999
Set-GoblinCell $ws 11 13 -5249859   # M11: -10249858 -> -5249859
Set-GoblinCell $ws 11 14 -5287   # N11: -5288 -> -5287
Set-GoblinCell $ws 61 8 4157.48   # H61: 3923.5925 -> 4157.48
Set-GoblinCell $ws 61 9 4226.5415   # I61: 3978.3462 -> 4226.5415
Set-GoblinCell $ws 61 11 4226.5415   # K61: 3978.3462 -> 4226.5415
Set-GoblinCell $ws 61 13 -4014.5415   # M61: -3766.3462 -> -4014.5415
Set-GoblinCell $ws 122 8 1650.4286   # H122: 1759 -> 1650.4286
Set-GoblinCell $ws 122 9 1650.4286   # I122: 1759 -> 1650.4286
Set-GoblinCell $ws 122 11 4951.2858   # K122: 5277 -> 4951.2858
Set-GoblinCell $ws 122 13 -2501.2858   # M122: -2827 -> -2501.2858
Set-GoblinCell $ws 132 8 2341.975   # H132: 2417.1538 -> 2341.975
Set-GoblinCell $ws 132 9 2182.838   # I132: 2259.861 -> 2182.838
Set-GoblinCell $ws 132 11 6548.514000000001   # K132: 6779.583 -> 6548.514000000001
Set-GoblinCell $ws 132 13 -4018.514000000001   # M132: -4249.583 -> -4018.514000000001
Set-GoblinCell $ws 136 8 4157.48   # H136: 3923.5925 -> 4157.48
Set-GoblinCell $ws 136 9 4226.5415   # I136: 3978.3462 -> 4226.5415
Set-GoblinCell $ws 136 11 12679.6245   # K136: 11935.0386 -> 12679.6245
Set-GoblinCell $ws 136 13 -10129.6245   # M136: -9385.0386 -> -10129.6245

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-GoblinCell $ws 22 8 908.2727   # H22: 688.1177 -> 908.2727
Set-GoblinCell $ws 22 9 956.125   # I22: 733.0833 -> 956.125
Set-GoblinCell $ws 22 10 780.6667   # J22: 580.2 -> 780.6667
Set-GoblinCell $ws 22 11 956.125   # K22: 733.0833 -> 956.125
Set-GoblinCell $ws 22 12 780.6667   # L22: 580.2 -> 780.6667
Set-GoblinCell $ws 22 13 -783.125   # M22: -560.0833 -> -783.125
Set-GoblinCell $ws 22 14 -1126.6667   # N22: -926.2 -> -1126.6667
Set-GoblinCell $ws 105 8 2547.7273   # H105: 2252.3572 -> 2547.7273
Set-GoblinCell $ws 105 9 1672.6666   # I105: 1505.5 -> 1672.6666
Set-GoblinCell $ws 105 10 3597.8   # J105: 3248.1667 -> 3597.8
Set-GoblinCell $ws 105 11 1672.6666   # K105: 1505.5 -> 1672.6666
Set-GoblinCell $ws 105 12 3597.8   # L105: 3248.1667 -> 3597.8
Set-GoblinCell $ws 105 13 74.33339999999998   # M105: 241.5 -> 74.33339999999998
Set-GoblinCell $ws 105 14 -7091.8   # N105: -6742.1667 -> -7091.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-GoblinCell $ws 19 8 1067.7142   # H19: 1075.625 -> 1067.7142
Set-GoblinCell $ws 19 9 412.33334   # I19: 474.75 -> 412.33334
Set-GoblinCell $ws 19 10 5000   # J19: 1676.5 -> 5000
Set-GoblinCell $ws 19 11 412.33334   # K19: 474.75 -> 412.33334
Set-GoblinCell $ws 19 12 5000   # L19: 1676.5 -> 5000
Set-GoblinCell $ws 19 13 -242.33334   # M19: -304.75 -> -242.33334
Set-GoblinCell $ws 19 14 -5340   # N19: -2016.5 -> -5340
Set-GoblinCell $ws 22 8 1363.7916   # H22: 1318.84 -> 1363.7916
Set-GoblinCell $ws 22 9 924.25   # I22: 884 -> 924.25
Set-GoblinCell $ws 22 11 924.25   # K22: 884 -> 924.25
Set-GoblinCell $ws 22 13 -574.25   # M22: -534 -> -574.25
Set-GoblinCell $ws 24 8 1067.7142   # H24: 1075.625 -> 1067.7142
Set-GoblinCell $ws 24 9 412.33334   # I24: 474.75 -> 412.33334
Set-GoblinCell $ws 24 10 5000   # J24: 1676.5 -> 5000
Set-GoblinCell $ws 24 11 412.33334   # K24: 474.75 -> 412.33334
Set-GoblinCell $ws 24 12 5000   # L24: 1676.5 -> 5000
Set-GoblinCell $ws 24 13 -242.33334   # M24: -304.75 -> -242.33334
Set-GoblinCell $ws 24 14 -5340   # N24: -2016.5 -> -5340
Set-GoblinCell $ws 94 8 2157.923   # H94: 1900.5 -> 2157.923
Set-GoblinCell $ws 94 9 2313.4   # I94: 2697 -> 2313.4
Set-GoblinCell $ws 94 10 2060.75   # J94: 1502.25 -> 2060.75
Set-GoblinCell $ws 94 11 2313.4   # K94: 2697 -> 2313.4
Set-GoblinCell $ws 94 12 2060.75   # L94: 1502.25 -> 2060.75
Set-GoblinCell $ws 94 13 -1862.4   # M94: -2246 -> -1862.4
Set-GoblinCell $ws 94 14 -2962.75   # N94: -2404.25 -> -2962.75
Set-GoblinCell $ws 99 8 2444.6   # H99: 2737.3333 -> 2444.6
Set-GoblinCell $ws 99 9 2444.6   # I99: 2737.3333 -> 2444.6
Set-GoblinCell $ws 99 11 2444.6   # K99: 2737.3333 -> 2444.6
Set-GoblinCell $ws 99 13 -946.5999999999999   # M99: -1239.3333 -> -946.5999999999999
Set-GoblinCell $ws 122 8 2037.4572   # H122: 2067.9707 -> 2037.4572
Set-GoblinCell $ws 122 9 1722.84   # I122: 1752.9584 -> 1722.84
Set-GoblinCell $ws 122 11 5168.52   # K122: 5258.8752 -> 5168.52
Set-GoblinCell $ws 122 13 -2718.52   # M122: -2808.8752 -> -2718.52
Set-GoblinCell $ws 126 8 2444.6   # H126: 2737.3333 -> 2444.6
Set-GoblinCell $ws 126 9 2444.6   # I126: 2737.3333 -> 2444.6
Set-GoblinCell $ws 126 11 7333.799999999999   # K126: 8211.999899999999 -> 7333.799999999999
Set-GoblinCell $ws 126 13 -4863.799999999999   # M126: -5741.999899999999 -> -4863.799999999999
Set-GoblinCell $ws 132 8 1917.36   # H132: 1956.0869 -> 1917.36
Set-GoblinCell $ws 132 9 1913.9166   # I132: 1956.0869 -> 1913.9166
Set-GoblinCell $ws 132 10 2000   # J132: 0 -> 2000
Set-GoblinCell $ws 132 11 5741.7498   # K132: 5868.2607 -> 5741.7498
Set-GoblinCell $ws 132 12 6000   # L132: 0 -> 6000
Set-GoblinCell $ws 132 13 -3211.7498   # M132: -3338.2607 -> -3211.7498
Set-GoblinCell $ws 132 14 -11060   # N132: None -> -11060

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-GoblinCell $ws 17 8 11222520   # H17: 10100268 -> 11222520
Set-GoblinCell $ws 17 10 333734   # J17: 250301.5 -> 333734
Set-GoblinCell $ws 17 12 1001202   # L17: 750904.5 -> 1001202
Set-GoblinCell $ws 17 14 -1001540   # N17: -751242.5 -> -1001540
Set-GoblinCell $ws 101 8 9666.583000000001   # H101: 9999.909 -> 9666.583000000001
Set-GoblinCell $ws 101 10 9666.583000000001   # J101: 9999.909 -> 9666.583000000001
Set-GoblinCell $ws 101 12 28999.749   # L101: 29999.727 -> 28999.749
Set-GoblinCell $ws 101 14 -33867.749   # N101: -34867.727 -> -33867.749
Set-GoblinCell $ws 128 8 202498   # H128: 0 -> 202498
Set-GoblinCell $ws 128 9 202498   # I128: 0 -> 202498
Set-GoblinCell $ws 128 11 607494   # K128: 0 -> 607494
Set-GoblinCell $ws 128 13 -602514   # M128: None -> -602514

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-GoblinCell $ws 107 8 359.1111   # H107: 391.875 -> 359.1111
Set-GoblinCell $ws 107 9 321.2857   # I107: 358.66666 -> 321.2857
Set-GoblinCell $ws 107 11 321.2857   # K107: 358.66666 -> 321.2857
Set-GoblinCell $ws 107 13 1598.7143   # M107: 1561.33334 -> 1598.7143
Set-GoblinCell $ws 122 8 6126.8213   # H122: 6181.1377 -> 6126.8213
Set-GoblinCell $ws 122 9 6632.75   # I122: 6701.4707 -> 6632.75
Set-GoblinCell $ws 122 10 5452.25   # J122: 5444 -> 5452.25
Set-GoblinCell $ws 122 11 19898.25   # K122: 20104.4121 -> 19898.25
Set-GoblinCell $ws 122 12 16356.75   # L122: 16332 -> 16356.75
Set-GoblinCell $ws 122 13 -17448.25   # M122: -17654.4121 -> -17448.25
Set-GoblinCell $ws 122 14 -21256.75   # N122: -21232 -> -21256.75
Set-GoblinCell $ws 126 8 3412.2856   # H126: 3519.75 -> 3412.2856
Set-GoblinCell $ws 126 9 2994.2   # I126: 3024 -> 2994.2
Set-GoblinCell $ws 126 10 4457.5   # J126: 5007 -> 4457.5
Set-GoblinCell $ws 126 11 8982.599999999999   # K126: 9072 -> 8982.599999999999
Set-GoblinCell $ws 126 12 13372.5   # L126: 15021 -> 13372.5
Set-GoblinCell $ws 126 13 -6512.599999999999   # M126: -6602 -> -6512.599999999999
Set-GoblinCell $ws 126 14 -18312.5   # N126: -19961 -> -18312.5
Set-GoblinCell $ws 132 8 2569.3784   # H132: 2632.4167 -> 2569.3784
Set-GoblinCell $ws 132 9 2408.6775   # I132: 2478.9666 -> 2408.6775
Set-GoblinCell $ws 132 11 7226.032499999999   # K132: 7436.899800000001 -> 7226.032499999999
Set-GoblinCell $ws 132 13 -4696.032499999999   # M132: -4906.899800000001 -> -4696.032499999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-GoblinCell $ws 2 8 999999   # H2: 0 -> 999999
Set-GoblinCell $ws 2 9 999999   # I2: 0 -> 999999
Set-GoblinCell $ws 2 10 999999   # J2: 0 -> 999999
Set-GoblinCell $ws 2 11 999999   # K2: 0 -> 999999
Set-GoblinCell $ws 2 12 999999   # L2: 0 -> 999999
Set-GoblinCell $ws 2 13 -999887   # M2: None -> -999887
Set-GoblinCell $ws 2 14 -1000223   # N2: None -> -1000223
Set-GoblinCell $ws 11 8 0   # H11: 4125 -> 0
Set-GoblinCell $ws 11 10 0   # J11: 4125 -> 0
Set-GoblinCell $ws 11 12 0   # L11: 4125 -> 0
Clear-GoblinCell $ws 11 14   # N11: -4405 -> (removed)
Set-GoblinCell $ws 13 8 11800   # H13: 8024.75 -> 11800
Set-GoblinCell $ws 13 10 20000   # J13: 9499.666999999999 -> 20000
Set-GoblinCell $ws 13 12 20000   # L13: 9499.666999999999 -> 20000
Set-GoblinCell $ws 13 14 -20280   # N13: -9779.666999999999 -> -20280
Set-GoblinCell $ws 19 8 97045.55   # H19: 106650.2 -> 97045.55
Set-GoblinCell $ws 19 9 116277.89   # I19: 130687.75 -> 116277.89
Set-GoblinCell $ws 19 11 116277.89   # K19: 130687.75 -> 116277.89
Set-GoblinCell $ws 19 13 -116107.89   # M19: -130517.75 -> -116107.89
Set-GoblinCell $ws 22 8 3302.2   # H22: 3458.111 -> 3302.2
Set-GoblinCell $ws 22 9 3352.7273   # I22: 3498 -> 3352.7273
Set-GoblinCell $ws 22 10 3272.9473   # J22: 3434.647 -> 3272.9473
Set-GoblinCell $ws 22 11 3352.7273   # K22: 3498 -> 3352.7273
Set-GoblinCell $ws 22 12 3272.9473   # L22: 3434.647 -> 3272.9473
Set-GoblinCell $ws 22 13 -3057.7273   # M22: -3203 -> -3057.7273
Set-GoblinCell $ws 22 14 -3862.9473   # N22: -4024.647 -> -3862.9473
Set-GoblinCell $ws 27 8 3302.2   # H27: 3458.111 -> 3302.2
Set-GoblinCell $ws 27 9 3352.7273   # I27: 3498 -> 3352.7273
Set-GoblinCell $ws 27 10 3272.9473   # J27: 3434.647 -> 3272.9473
Set-GoblinCell $ws 27 11 3352.7273   # K27: 3498 -> 3352.7273
Set-GoblinCell $ws 27 12 3272.9473   # L27: 3434.647 -> 3272.9473
Set-GoblinCell $ws 27 13 -3245.7273   # M27: -3391 -> -3245.7273
Set-GoblinCell $ws 27 14 -3486.9473   # N27: -3648.647 -> -3486.9473
Set-GoblinCell $ws 46 8 2399.625   # H46: 2357.2122 -> 2399.625
Set-GoblinCell $ws 46 10 3158.8948   # J46: 3050.95 -> 3158.8948
Set-GoblinCell $ws 46 12 3158.8948   # L46: 3050.95 -> 3158.8948
Set-GoblinCell $ws 46 14 -3534.8948   # N46: -3426.95 -> -3534.8948
Set-GoblinCell $ws 55 8 704.46155   # H55: 685.8570999999999 -> 704.46155
Set-GoblinCell $ws 55 9 347   # I55: 337.36365 -> 347
Set-GoblinCell $ws 55 10 1896   # J55: 1963.6666 -> 1896
Set-GoblinCell $ws 55 11 347   # K55: 337.36365 -> 347
Set-GoblinCell $ws 55 12 1896   # L55: 1963.6666 -> 1896
Set-GoblinCell $ws 55 13 -174   # M55: -164.36365 -> -174
Set-GoblinCell $ws 55 14 -2242   # N55: -2309.6666 -> -2242
Set-GoblinCell $ws 93 8 2307.6924   # H93: 2364.3948 -> 2307.6924
Set-GoblinCell $ws 93 9 1605.8846   # I93: 1664 -> 1605.8846
Set-GoblinCell $ws 93 11 1605.8846   # K93: 1664 -> 1605.8846
Set-GoblinCell $ws 93 13 -357.8846000000001   # M93: -416 -> -357.8846000000001
Set-GoblinCell $ws 122 8 5658.6665   # H122: 6257.8335 -> 5658.6665
Set-GoblinCell $ws 122 9 5406.75   # I122: 5709.5 -> 5406.75
Set-GoblinCell $ws 122 10 6666.3335   # J122: 8999.5 -> 6666.3335
Set-GoblinCell $ws 122 11 16220.25   # K122: 17128.5 -> 16220.25
Set-GoblinCell $ws 122 12 19999.0005   # L122: 26998.5 -> 19999.0005
Set-GoblinCell $ws 122 13 -13770.25   # M122: -14678.5 -> -13770.25
Set-GoblinCell $ws 122 14 -24899.0005   # N122: -31898.5 -> -24899.0005
Set-GoblinCell $ws 132 8 3624.4   # H132: 3773.5264 -> 3624.4
Set-GoblinCell $ws 132 9 3489.8635   # I132: 3759.75 -> 3489.8635
Set-GoblinCell $ws 132 11 10469.5905   # K132: 11279.25 -> 10469.5905
Set-GoblinCell $ws 132 13 -7939.5905   # M132: -8749.25 -> -7939.5905
Set-GoblinCell $ws 136 8 4756.68   # H136: 4947.617 -> 4756.68
Set-GoblinCell $ws 136 9 4531.8   # I136: 4918.8887 -> 4531.8
Set-GoblinCell $ws 136 10 4906.6   # J136: 4965.448 -> 4906.6
Set-GoblinCell $ws 136 11 13595.4   # K136: 14756.6661 -> 13595.4
Set-GoblinCell $ws 136 12 14719.8   # L136: 14896.344 -> 14719.8
Set-GoblinCell $ws 136 13 -11045.4   # M136: -12206.6661 -> -11045.4
Set-GoblinCell $ws 136 14 -19819.8   # N136: -19996.344 -> -19819.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-GoblinCell $ws 24 8 0   # H24: 9000 -> 0
Set-GoblinCell $ws 24 9 0   # I24: 9000 -> 0
Set-GoblinCell $ws 24 11 0   # K24: 9000 -> 0
Clear-GoblinCell $ws 24 13   # M24: -8770 -> (removed)
Set-GoblinCell $ws 132 8 3153.8125   # H132: 3237.3547 -> 3153.8125
Set-GoblinCell $ws 132 9 2497.4285   # I132: 2569.037 -> 2497.4285
Set-GoblinCell $ws 132 11 7492.2855   # K132: 7707.110999999999 -> 7492.2855
Set-GoblinCell $ws 132 13 -4962.2855   # M132: -5177.110999999999 -> -4962.2855

Write-Host "Applied all Goblin_Profits cell updates."